$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8, 9, 10: fill in the "N°" (order number) column A ---
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

# --- Insert a new row at 11 (push everything from 11 down to 12, including
#     the old last-border row which ends up on row 22). Excel's default
#     shift-down copies the style from the row above (row 10), so fix the
#     new row 11's style by copying it from row 12 (the old row 11, which
#     kept the regular body-row style). ---
$ws.Rows("11:11").Insert(-4121)  # xlShiftDown
$ws.Range("A12:E12").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 11 : "Contrôle du contenu du champs email" ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Contrôle du contenu du champs email"
$ws.Range("C11").Value = "Récupération du contenu du champs email et vérification de son contenu par le biaix d'un REGEX"
$ws.Range("D11").Value = "Affichage d'une alerte si le format de l'adresse email ne correspond pas au format attendu"
$ws.Range("E11").Value = "OK / Affichage d'une alerte indiquant : `"Le mail saisie ne semble pas être valide. Merci de saisir une adresse mail valide.`" "
$ws.Rows(11).RowHeight = 87

# --- Row 12 : "Commander" (content replaces the former row-11 content) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Commander"
$ws.Range("C12").Value = "Récupération du contenu des différents champs du formulaires, contrôle de la validité du contenu du champs email, envoi de la requete au serveur, récupération du N° de commande  et ouverture de la page confirmation"
$ws.Range("D12").Value = "Ouverture de la page confirmation"
$ws.Range("E12").Value = "OK"
$ws.Rows(12).RowHeight = 130.5

# --- Row 13 : "Affichage de la confirmation de commande" (former row-12 content) ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Affichage de la confirmation de commande"
$ws.Range("C13").Value = "Récupération du numéro de commande dans l'adresse de la page et affichage du numéro de commande"
$ws.Range("D13").Value = "Affichage du numéro de commande"
$ws.Range("E13").Value = "OK"
$ws.Rows(13).RowHeight = 65.25

# --- View state: frozen pane / active selection ---
$ws.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D18").Select()

$wb.Save()
